# Generate Report for Archive
#
# The localization status for these rows moved from "Ready for handoff" to
# "In Translation". That text lives in the "Status" column, which is shown
# as the zh-cn/de-de columns on the Overview sheet and as the "Status"
# column on each per-language sheet. Updating the text causes the host
# report generator to re-fit those (now narrower) columns.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

# --- zh-cn sheet: column C is "Status" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: column C is "Status" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
